# Recreated mantel correlograms with Euclidean distances:
# update the "Mantel r" and "p" columns of the correlogram table
# (table 1, columns 3 and 4) for each distance-class row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 2; Col = 3; Value = "-0.022" },
    @{ Row = 2; Col = 4; Value = "0.247" },
    @{ Row = 3; Col = 3; Value = "0.052" },
    @{ Row = 3; Col = 4; Value = "0.058" },
    @{ Row = 4; Col = 3; Value = "-0.007" },
    @{ Row = 4; Col = 4; Value = "0.494" },
    @{ Row = 5; Col = 3; Value = "-0.042" },
    @{ Row = 5; Col = 4; Value = "0.237" },
    @{ Row = 6; Col = 3; Value = "-0.050" },
    @{ Row = 6; Col = 4; Value = "0.2" },
    @{ Row = 7; Col = 3; Value = "-0.003" },
    @{ Row = 7; Col = 4; Value = "0.773" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Value
}
